# Auto-generated edit script applying cached-value updates to the
# "Garuda_Profits" leve-profit workbook (FFXIV leve profitability sheets).
# Each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) holds refreshed
# market-board price snapshots in columns H:N for specific leve rows.
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 94.111115
$ws.Range("I2").Value = 94.111115
$ws.Range("K2").Value = 94.111115
$ws.Range("M2").Value = 18.888885

$ws.Range("H21").Value = 22200
$ws.Range("I21").Value = 14333.333
$ws.Range("K21").Value = 14333.333
$ws.Range("M21").Value = -13865.333

$ws.Range("H23").Value = 22200
$ws.Range("I23").Value = 14333.333
$ws.Range("K23").Value = 14333.333
$ws.Range("M23").Value = -14099.333

$ws.Range("H43").Value = 83333736
$ws.Range("I43").Value = 111111480
$ws.Range("J43").Value = 33333794
$ws.Range("K43").Value = 111111480
$ws.Range("L43").Value = 33333794
$ws.Range("M43").Value = -111111411
$ws.Range("N43").Value = -33333932

$ws.Range("H87").Value = 34100
$ws.Range("J87").Value = 33116.668
$ws.Range("L87").Value = 33116.668
$ws.Range("N87").Value = -35612.668

$ws.Range("H90").Value = 34100
$ws.Range("J90").Value = 33116.668
$ws.Range("L90").Value = 99350.00399999999
$ws.Range("N90").Value = -111830.004

$ws.Range("H98").Value = 17288
$ws.Range("I98").Value = 19669.334
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 19669.334
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = -18171.334
$ws.Range("N98").Value = -5996

$ws.Range("H112").Value = 1618.6875
$ws.Range("J112").Value = 1653.2667
$ws.Range("L112").Value = 4959.800099999999
$ws.Range("N112").Value = -7175.800099999999

$ws.Range("H122").Value = 17288
$ws.Range("I122").Value = 19669.334
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 59008.00199999999
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -56558.00199999999
$ws.Range("N122").Value = -13900

$ws.Range("H131").Value = 961.5238000000001
$ws.Range("I131").Value = 886.13336
$ws.Range("J131").Value = 1150
$ws.Range("K131").Value = 2658.40008
$ws.Range("L131").Value = 3450
$ws.Range("M131").Value = 2381.59992
$ws.Range("N131").Value = -13530

$ws.Range("H137").Value = 1236.7858
$ws.Range("I137").Value = 965.9677
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 2897.9031
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -347.9031
$ws.Range("N137").Value = -11100


# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18945.719
$ws.Range("I32").Value = 18726.152
$ws.Range("J32").Value = 21580.5
$ws.Range("K32").Value = 18726.152
$ws.Range("L32").Value = 21580.5
$ws.Range("M32").Value = -18439.152
$ws.Range("N32").Value = -22154.5

$ws.Range("H61").Value = 1723.973
$ws.Range("I61").Value = 1432.9615
$ws.Range("J61").Value = 2411.818
$ws.Range("K61").Value = 1432.9615
$ws.Range("L61").Value = 2411.818
$ws.Range("M61").Value = -1220.9615
$ws.Range("N61").Value = -2835.818

$ws.Range("H74").Value = 806.7213
$ws.Range("I74").Value = 833
$ws.Range("J74").Value = 654.8889
$ws.Range("K74").Value = 833
$ws.Range("L74").Value = 654.8889
$ws.Range("M74").Value = 41
$ws.Range("N74").Value = -2402.8889

$ws.Range("H77").Value = 806.7213
$ws.Range("I77").Value = 833
$ws.Range("J77").Value = 654.8889
$ws.Range("K77").Value = 4165
$ws.Range("L77").Value = 3274.4445
$ws.Range("M77").Value = 203
$ws.Range("N77").Value = -12010.4445

$ws.Range("H102").Value = 1328.25
$ws.Range("I102").Value = 1282.2222
$ws.Range("J102").Value = 1466.3334
$ws.Range("K102").Value = 1282.2222
$ws.Range("L102").Value = 1466.3334
$ws.Range("M102").Value = 339.7778000000001
$ws.Range("N102").Value = -4710.3334

$ws.Range("H132").Value = 6021.6875
$ws.Range("I132").Value = 6582.41
$ws.Range("J132").Value = 3591.889
$ws.Range("K132").Value = 19747.23
$ws.Range("L132").Value = 10775.667
$ws.Range("M132").Value = -17217.23
$ws.Range("N132").Value = -15835.667

$ws.Range("H136").Value = 1723.973
$ws.Range("I136").Value = 1432.9615
$ws.Range("J136").Value = 2411.818
$ws.Range("K136").Value = 4298.8845
$ws.Range("L136").Value = 7235.454000000001
$ws.Range("M136").Value = -1748.8845
$ws.Range("N136").Value = -12335.454


# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4524.079
$ws.Range("I20").Value = 5132.68
$ws.Range("J20").Value = 3353.6924
$ws.Range("K20").Value = 5132.68
$ws.Range("L20").Value = 3353.6924
$ws.Range("M20").Value = -4885.68
$ws.Range("N20").Value = -3847.6924

$ws.Range("H64").Value = 479.54544
$ws.Range("I64").Value = 472
$ws.Range("J64").Value = 483.85715
$ws.Range("K64").Value = 472
$ws.Range("L64").Value = 483.85715
$ws.Range("M64").Value = -247
$ws.Range("N64").Value = -933.85715

$ws.Range("H67").Value = 479.54544
$ws.Range("I67").Value = 472
$ws.Range("J67").Value = 483.85715
$ws.Range("K67").Value = 472
$ws.Range("L67").Value = 483.85715
$ws.Range("M67").Value = 308
$ws.Range("N67").Value = -2043.85715

$ws.Range("H134").Value = 4398.9067
$ws.Range("I134").Value = 5034.2905
$ws.Range("J134").Value = 2757.5
$ws.Range("K134").Value = 15102.8715
$ws.Range("L134").Value = 8272.5
$ws.Range("M134").Value = -12567.8715
$ws.Range("N134").Value = -13342.5


# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 37149.75
$ws.Range("J68").Value = 39533
$ws.Range("L68").Value = 39533
$ws.Range("N68").Value = -41031

$ws.Range("H71").Value = 37149.75
$ws.Range("J71").Value = 39533
$ws.Range("L71").Value = 118599
$ws.Range("N71").Value = -126087

$ws.Range("H74").Value = 24999.5
$ws.Range("J74").Value = 39999
$ws.Range("L74").Value = 39999
$ws.Range("N74").Value = -41747

$ws.Range("H77").Value = 24999.5
$ws.Range("J77").Value = 39999
$ws.Range("L77").Value = 119997
$ws.Range("N77").Value = -128733


# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3120
$ws.Range("I80").Value = 900
$ws.Range("J80").Value = 3675
$ws.Range("K80").Value = 2700
$ws.Range("L80").Value = 11025
$ws.Range("M80").Value = -1764
$ws.Range("N80").Value = -12897

$ws.Range("H83").Value = 3120
$ws.Range("I83").Value = 900
$ws.Range("J83").Value = 3675
$ws.Range("K83").Value = 8100
$ws.Range("L83").Value = 33075
$ws.Range("M83").Value = -3420
$ws.Range("N83").Value = -42435

$ws.Range("H129").Value = 8773622
$ws.Range("I129").Value = 840.7778
$ws.Range("J129").Value = 16669125
$ws.Range("K129").Value = 2522.3334
$ws.Range("L129").Value = 50007375
$ws.Range("M129").Value = 2477.6666
$ws.Range("N129").Value = -50017375

$ws.Range("H131").Value = 594.42426
$ws.Range("J131").Value = 803.9508
$ws.Range("L131").Value = 2411.8524
$ws.Range("N131").Value = -12491.8524

$ws.Range("H132").Value = 1798.3871
$ws.Range("I132").Value = 900
$ws.Range("J132").Value = 2226.1904
$ws.Range("K132").Value = 8100
$ws.Range("L132").Value = 20035.7136
$ws.Range("M132").Value = -5570
$ws.Range("N132").Value = -25095.7136


# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2033.5714
$ws.Range("I97").Value = 2197.2727
$ws.Range("J97").Value = 1433.3334
$ws.Range("K97").Value = 2197.2727
$ws.Range("L97").Value = 1433.3334
$ws.Range("M97").Value = -1701.2727
$ws.Range("N97").Value = -2425.3334

$ws.Range("H122").Value = 30305282
$ws.Range("I122").Value = 40002012
$ws.Range("K122").Value = 120006036
$ws.Range("M122").Value = -120003586


# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5212.0605
$ws.Range("I136").Value = 5731.077
$ws.Range("K136").Value = 17193.231
$ws.Range("M136").Value = -14643.231


# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 125000450
$ws.Range("I96").Value = 125000450
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 125000450
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -124999077
$ws.Range("N96").ClearContents()

$ws.Range("H136").Value = 1307.452
$ws.Range("I136").Value = 1136.7646
$ws.Range("J136").Value = 1703.1364
$ws.Range("K136").Value = 3410.2938
$ws.Range("L136").Value = 5109.4092
$ws.Range("M136").Value = -860.2937999999999
$ws.Range("N136").Value = -10209.4092
